# Sacramento Seminarians workbook — add a new member row.
#
# The "members" sheet is sorted in the order entries were added, and a new
# member ("Leo" / last name, "Petroni" / first name) needs to be inserted
# just above the "Purseglove" row (which currently sits at row 69), pushing
# every row below it down by one. Only the last-name and first-name columns
# are known for this entry; email/phone stay blank, same as several other
# rows in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("members")

# Shift row 69 (and everything after it) down by inserting a fresh row;
# Excel copies the formatting of the row above into the new row.
$ws.Rows.Item(69).Insert()

$ws.Range("A69").Value = "Leo"
$ws.Range("B69").Value = "Petroni"

# Match the row height used throughout the rest of the table.
$ws.Rows.Item(69).RowHeight = 18.75

# No email/phone for this entry.
$ws.Range("C69").ClearContents()
$ws.Range("D69").ClearContents()
$ws.Range("E69").ClearContents()
